$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing GDP values for rows 12-30 (column B)
$ws.Range("B12").Value = 227080.47500000001
$ws.Range("B13").Value = 236415.902
$ws.Range("B14").Value = 243606.20699999999
$ws.Range("B15").Value = 263007.68400000001
$ws.Range("B16").Value = 284629.81400000001
$ws.Range("B17").Value = 314651.88199999998
$ws.Range("B18").Value = 332946.79300000001
$ws.Range("B19").Value = 343969.79700000002
$ws.Range("B20").Value = 327718.81599999999
$ws.Range("B21").Value = 340179.679
$ws.Range("B22").Value = 359386.45799999998
$ws.Range("B23").Value = 377846.40700000001
$ws.Range("B24").Value = 396428.07199999999
$ws.Range("B25").Value = 417730.68800000002
$ws.Range("B26").Value = 438084.24800000002
$ws.Range("B27").Value = 456712.67499999999
$ws.Range("B28").Value = 480906.20299999998
$ws.Range("B29").Value = 511962.53700000001
$ws.Range("B30").Value = 540633.96100000001

# Add new row 31 with 2020 data point
$ws.Range("A31").Value = 43831
$ws.Range("A31").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B31").Value = 534807.25100000005
$ws.Range("B31").NumberFormat = "0.000"

# Update selection to match target state (user selected full columns A:B).
# The COM layer always anchors the active cell to the top-left corner of the
# selected rectangle, so we select the full A:B columns (A1:B1048576) to
# match the saved sqref as closely as possible.
$ws.Range("A1:B1048576").Select()
